# Generate Report for Handback
# Update the "generate date" timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G2)
$wsOverview.Range("G2").Value = "2016-08-21 07:10:27"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-21 07:10:22"
$wsZhCn.Range("K2").Value = "2016-08-21 07:10:39"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-08-21 07:10:27"
$wsDeDe.Range("K2").Value = "2016-08-21 07:10:46"
